$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# --- Rename the "MTRN1" test group to "MTRN1.1" for the existing rows (38-40) ---
$ws.Range("B38").Value = "MTRN1.1"
$ws.Range("B39").Value = "MTRN1.1"
$ws.Range("B40").Value = "MTRN1.1"

# --- Fill in the logfile name for rows 38-40, now that it is known ---
$ws.Range("H38").Value = "09.35.txt"
$ws.Range("H39").Value = "09.35.txt"
$ws.Range("H40").Value = "09.35.txt"

# --- Add three new rows (41-43) for the MTRN1.2 re-run (bug-fix) results ---
$ws.Range("A41").Value = "MTRN1.2.1"
$ws.Range("A42").Value = "MTRN1.2.2"
$ws.Range("A43").Value = "MTRN1.2.3"
$ws.Range("B41:B43").Value = "MTRN1.2"

$d = Get-Date -Year 2016 -Month 6 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("G41").Value = $d
$ws.Range("G42").Value = $d
$ws.Range("G43").Value = $d

# Extend the shared formulas down through the new rows
$ws.Range("D2:D43").Formula = "=1/3"
$ws.Range("E41:E43").Formula = "=C41/D41"
$ws.Range("F41:F43").Formula = "=AVERAGEIF(B:B,B41,E:E)"

# Match the number formatting used by the neighboring rows
$ws.Range("E41:F43").NumberFormat = "0.00"

$excel.CalculateFullRebuild()

# --- Conditional formatting bookkeeping: the workbook's dxf table grows by two
#     additional "good" (green) style entries, matching the existing top10 rules ---
for ($i = 0; $i -lt 2; $i++) {
    $tmp = $ws.Range("Z1:Z2")
    $fc = $tmp.FormatConditions.AddTop10()
    $fc.TopBottom = 1
    $fc.Rank = 10
    $fc.Percent = $true
    $fc.Font.Color = 24832
    $fc.Interior.Color = 13561798
    $fc.Delete()
}

$ws.Application.Goto($ws.Range("H42"))
$ws.Range("H42").Select()
